$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.432.16'
$ws.Range("E2").Value = '  +1.65%  '

$ws.Range("D3").Value = '1.862.93'
$ws.Range("E3").Value = '  +0.76%  '

$ws.Range("D4").Formula = '="1.012"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").Formula = '="311.41"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +0.65%  '

$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("D7").Formula = '="0.4773"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  +3.32%  '

$ws.Range("D9").Formula = '="0.07324"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  +1.41%  '

$ws.Range("D10").Formula = '="0.9335"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  +0.41%  '

$ws.Range("E11").Value = '  +5.11%  '

$ws.Range("D12").Formula = '="0.07803"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  +0.74%  '

$ws.Range("D13").Value = '1.881.16'
$ws.Range("E13").Value = '  +2.39%  '

$ws.Range("D14").Formula = '="5.434"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  +1.71%  '

$ws.Range("D15").Formula = '="6.550"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  +1.69%  '

$ws.Range("E16").Value = '  +1.71%  '

$ws.Range("E17").Value = '  -0.20%  '

$ws.Range("D18").Formula = '="0.000008803"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  +1.64%  '

$ws.Range("E19").Value = '  -0.18%  '

$ws.Range("D20").Value = '27.529.57'
$ws.Range("E20").Value = '  +1.92%  '

$ws.Range("D21").Formula = '="14.63"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  +1.15%  '

$ws.Range("D22").Formula = '="5.115"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)

$ws.Range("E23").Value = '  +0.48%  '

$ws.Range("D24").Formula = '="1.941"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +0.63%  '

$ws.Range("D25").Formula = '="155.81"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +1.91%  '

$ws.Range("D26").Formula = '="18.47"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  +1.27%  '

$ws.Range("D27").Formula = '="2.019"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +0.72%  '

$ws.Range("D28").Formula = '="115.30"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  +0.94%  '

$ws.Range("D29").Formula = '="4.942"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  -0.28%  '

$ws.Range("D30").Formula = '="0.08878"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  +0.03%  '

$ws.Range("D31").Formula = '="3.327"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  +0.06%  '

$ws.Range("E32").Value = '  +3.44%  '

$ws.Range("D33").Formula = '="0.7580"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  +2.10%  '

$ws.Range("D34").Formula = '="4.598"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)

$ws.Range("D35").Formula = '="2.720"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -0.99%  '

$ws.Range("D36").Formula = '="0.02047"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  +4.14%  '

$ws.Range("E37").Value = '  +0.69%  '

$ws.Range("D38").Formula = '="0.5561"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  +6.68%  '

$ws.Range("D39").Formula = '="0.05267"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  -0.07%  '

$ws.Range("D40").Formula = '="2.988"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  +0.25%  '

$ws.Range("D41").Formula = '="7.052"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  +0.85%  '

$ws.Range("D42").Formula = '="8.631"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  +4.66%  '

$ws.Range("E43").Value = '  +0.72%  '

$ws.Range("D44").Formula = '="0.4888"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +3.08%  '

$ws.Range("D45").Formula = '="10.72"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +0.75%  '

$ws.Range("E46").Value = '  -0.18%  '

$ws.Range("D47").Formula = '="103.05"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  +1.17%  '

$ws.Range("D48").Formula = '="1.655"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +2.83%  '

$ws.Range("D49").Formula = '="67.45"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  +2.52%  '

$ws.Range("E50").Value = '  +0.29%  '

$ws.Range("D51").Formula = '="0.9156"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  +3.13%  '

$excel.CutCopyMode = $false